# Auto-generated cell updates for Jogos_da_Semana_FlashScore_2025-02-14.xlsx
# Applies updated odds values per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("N2").Value = 15
$ws.Range("S2").Value = 1.88
$ws.Range("T2").Value = 2.02

# Row 3
$ws.Range("G3").Value = 2.1
$ws.Range("H3").Value = 3.3
$ws.Range("I3").Value = 3.7
$ws.Range("K3").Value = 2.1
$ws.Range("O3").Value = 1.3
$ws.Range("P3").Value = 3.5
$ws.Range("Q3").Value = 2.04
$ws.Range("R3").Value = 1.86
$ws.Range("AA3").Value = 7.5
$ws.Range("AG3").Value = 9.5
$ws.Range("AH3").Value = 6
$ws.Range("AL3").Value = 11

# Row 5
$ws.Range("G5").Value = 1.75
$ws.Range("I5").Value = 5.25
$ws.Range("Y5").Value = 2.2
$ws.Range("Z5").Value = 1.62
$ws.Range("AB5").Value = 7
$ws.Range("AI5").Value = 21
$ws.Range("AJ5").Value = 81
$ws.Range("AL5").Value = 11
$ws.Range("AR5").Value = 1.82
$ws.Range("AS5").Value = 2.08

# Row 6
$ws.Range("Q6").Value = 1.73
$ws.Range("R6").Value = 2.08

# Row 7
$ws.Range("Q7").Value = 2.63
$ws.Range("R7").Value = 1.5
$ws.Range("AG7").Value = 6.5
$ws.Range("AO7").Value = 51
$ws.Range("AR7").Value = 1.95
$ws.Range("AS7").Value = 1.95

# Row 8
$ws.Range("L8").Value = 3.25
$ws.Range("M8").Value = 1.11
$ws.Range("N8").Value = 6.5
$ws.Range("W8").Value = 1.62
$ws.Range("X8").Value = 2.2
$ws.Range("AB8").Value = 13
$ws.Range("AO8").Value = 23

# Row 10
$ws.Range("N10").Value = 5
$ws.Range("W10").Value = 1.75
$ws.Range("X10").Value = 2.05

# Row 11
$ws.Range("G11").Value = 1.44
$ws.Range("I11").Value = 7
$ws.Range("S11").Value = 1.95
$ws.Range("T11").Value = 1.9
$ws.Range("U11").Value = 2.38
$ws.Range("V11").Value = 1.57
$ws.Range("AA11").Value = 9
$ws.Range("AD11").Value = 10
$ws.Range("AG11").Value = 17
$ws.Range("AH11").Value = 9
$ws.Range("AI11").Value = 17

# Row 12
$ws.Range("G12").Value = 5
$ws.Range("H12").Value = 3.75
$ws.Range("I12").Value = 1.7
$ws.Range("K12").Value = 2.25
$ws.Range("W12").Value = 1.36
$ws.Range("X12").Value = 3
$ws.Range("AI12").Value = 15
$ws.Range("AL12").Value = 7.5

# Row 13
$ws.Range("G13").Value = 4.1
$ws.Range("H13").Value = 3.8
$ws.Range("I13").Value = 1.7
$ws.Range("J13").Value = 4.5
$ws.Range("N13").Value = 17
$ws.Range("O13").Value = 1.17
$ws.Range("P13").Value = 5
$ws.Range("Q13").Value = 1.57
$ws.Range("R13").Value = 2.35
$ws.Range("S13").Value = 1.95
$ws.Range("T13").Value = 1.9
$ws.Range("U13").Value = 2.38
$ws.Range("V13").Value = 1.53
$ws.Range("AA13").Value = 17
$ws.Range("AB13").Value = 26
$ws.Range("AC13").Value = 15
$ws.Range("AF13").Value = 34

# Row 14
$ws.Range("AG14").Value = 6.5
$ws.Range("AI14").Value = 17
$ws.Range("AL14").Value = 7.5
$ws.Range("AP14").Value = 29

# Row 15
$ws.Range("G15").Value = 2.1
$ws.Range("H15").Value = 3.5
$ws.Range("I15").Value = 3.2
$ws.Range("J15").Value = 2.75
$ws.Range("L15").Value = 3.75
$ws.Range("T15").Value = 1.53
$ws.Range("Y15").Value = 1.73
$ws.Range("Z15").Value = 2
$ws.Range("AA15").Value = 8
$ws.Range("AB15").Value = 10
$ws.Range("AD15").Value = 19
$ws.Range("AE15").Value = 17
$ws.Range("AH15").Value = 6.5
$ws.Range("AJ15").Value = 41
$ws.Range("AK15").Value = 201
$ws.Range("AL15").Value = 10
$ws.Range("AM15").Value = 17
$ws.Range("AO15").Value = 34
$ws.Range("AP15").Value = 26

# Row 18
$ws.Range("G18").Value = 2.9
$ws.Range("H18").Value = 3.1
$ws.Range("I18").Value = 2.38
$ws.Range("J18").Value = 3.6
$ws.Range("L18").Value = 3.2
$ws.Range("O18").Value = 1.36
$ws.Range("P18").Value = 3
$ws.Range("Q18").Value = 2.15
$ws.Range("R18").Value = 1.67
$ws.Range("AA18").Value = 8.5
$ws.Range("AD18").Value = 29
$ws.Range("AG18").Value = 8.5
$ws.Range("AK18").Value = 301
$ws.Range("AN18").Value = 10
$ws.Range("AO18").Value = 23
$ws.Range("AP18").Value = 21
$ws.Range("AQ18").Value = 34

# Row 20
$ws.Range("G20").Value = 1.17
$ws.Range("H20").Value = 6.2
$ws.Range("J20").Value = 1.55
$ws.Range("K20").Value = 2.6
$ws.Range("L20").Value = 12
$ws.Range("M20").Value = 1.04
$ws.Range("N20").Value = 9.5
$ws.Range("O20").Value = 1.19
$ws.Range("P20").Value = 4.35
$ws.Range("Q20").Value = 1.57
$ws.Range("R20").Value = 2.3
$ws.Range("U20").Value = 2.4
$ws.Range("V20").Value = 1.53
$ws.Range("W20").Value = 1.31
$ws.Range("X20").Value = 3.25
$ws.Range("Y20").Value = 2.42
$ws.Range("Z20").Value = 1.5
$ws.Range("AA20").Value = 6
$ws.Range("AB20").Value = 5.7
$ws.Range("AD20").Value = 6.5
$ws.Range("AE20").Value = 12.5
$ws.Range("AF20").Value = 45
$ws.Range("AG20").Value = 9.5
$ws.Range("AH20").Value = 14
$ws.Range("AI20").Value = 40
$ws.Range("AJ20").Value = 250
$ws.Range("AL20").Value = 32
$ws.Range("AP20").Value = 400

# Row 22
$ws.Range("G22").Value = 5
$ws.Range("H22").Value = 4
$ws.Range("L22").Value = 2.3
$ws.Range("M22").Value = 1.06
$ws.Range("N22").Value = 9.5
$ws.Range("Q22").Value = 2.05
$ws.Range("R22").Value = 1.75
$ws.Range("W22").Value = 1.44
$ws.Range("X22").Value = 2.63
$ws.Range("Y22").Value = 2
$ws.Range("Z22").Value = 1.73
$ws.Range("AC22").Value = 17
$ws.Range("AG22").Value = 9.5
$ws.Range("AH22").Value = 7.5
$ws.Range("AJ22").Value = 67
$ws.Range("AL22").Value = 6
$ws.Range("AM22").Value = 7
$ws.Range("AP22").Value = 15

# Row 23
$ws.Range("J23").Value = 3
$ws.Range("K23").Value = 1.95
$ws.Range("AG23").Value = 7.5

# Row 26
$ws.Range("N26").Value = 8
$ws.Range("W26").Value = 1.53
$ws.Range("X26").Value = 2.38
$ws.Range("AG26").Value = 7

# Row 27
$ws.Range("O27").Value = 1.18
$ws.Range("P27").Value = 4.5
$ws.Range("Q27").Value = 1.67
$ws.Range("R27").Value = 2.15
$ws.Range("U27").Value = 2.5
$ws.Range("V27").Value = 1.5

# Row 29
$ws.Range("M29").Value = 1.08
$ws.Range("N29").Value = 8
$ws.Range("Q29").Value = 2.25
$ws.Range("R29").Value = 1.62
$ws.Range("U29").Value = 4
$ws.Range("V29").Value = 1.22

# Row 32
$ws.Range("G32").Value = 2.55
$ws.Range("I32").Value = 2.63
$ws.Range("J32").Value = 3.25
$ws.Range("AD32").Value = 26
$ws.Range("AL32").Value = 8.5
$ws.Range("AN32").Value = 10
$ws.Range("AO32").Value = 26
$ws.Range("AP32").Value = 21

# Row 33
$ws.Range("AD33").Value = 21

# Row 34
$ws.Range("W34").Value = 1.29
$ws.Range("X34").Value = 3.5
$ws.Range("Y34").Value = 1.57
$ws.Range("Z34").Value = 2.25
$ws.Range("AE34").Value = 12
$ws.Range("AK34").Value = 126

# Row 35
$ws.Range("G35").Value = 1.91
$ws.Range("H35").Value = 3.25
$ws.Range("I35").Value = 3.85
$ws.Range("J35").Value = 2.47
$ws.Range("K35").Value = 2.12
$ws.Range("L35").Value = 4.15
$ws.Range("O35").Value = 1.31
$ws.Range("P35").Value = 2.9
$ws.Range("Q35").Value = 1.91
$ws.Range("U35").Value = 3.1
$ws.Range("V35").Value = 1.27
$ws.Range("Z35").Value = 1.85
$ws.Range("AA35").Value = 6.7
$ws.Range("AB35").Value = 8.75
$ws.Range("AC35").Value = 8.5
$ws.Range("AD35").Value = 16.5
$ws.Range("AE35").Value = 16
$ws.Range("AF35").Value = 28
$ws.Range("AH35").Value = 6.3
$ws.Range("AI35").Value = 14.5
$ws.Range("AK35").Value = 500
$ws.Range("AL35").Value = 11
$ws.Range("AM35").Value = 22
$ws.Range("AN35").Value = 12.5
$ws.Range("AO35").Value = 60
$ws.Range("AP35").Value = 35

# Row 37
$ws.Range("G37").Value = 3.4
$ws.Range("H37").Value = 3.25
$ws.Range("I37").Value = 1.95
$ws.Range("AI37").Value = 17
$ws.Range("AK37").Value = 401
$ws.Range("AM37").Value = 9

# Row 38
$ws.Range("M38").Value = 1.1
$ws.Range("N38").Value = 7

# Row 39
$ws.Range("G39").Value = 1.38
$ws.Range("H39").Value = 4.5
$ws.Range("I39").Value = 5.5
$ws.Range("J39").Value = 1.95
$ws.Range("K39").Value = 2.3
$ws.Range("L39").Value = 7.5
$ws.Range("Q39").Value = 1.93
$ws.Range("R39").Value = 1.93
$ws.Range("Y39").Value = 2.2
$ws.Range("Z39").Value = 1.62
$ws.Range("AB39").Value = 6
$ws.Range("AD39").Value = 9
$ws.Range("AH39").Value = 9.5
$ws.Range("AL39").Value = 13
$ws.Range("AM39").Value = 34
$ws.Range("AN39").Value = 19
$ws.Range("AO39").Value = 81

# Row 40
$ws.Range("J40").Value = 2.88
$ws.Range("M40").Value = 1.08
$ws.Range("N40").Value = 8
$ws.Range("AB40").Value = 9
$ws.Range("AF40").Value = 34
$ws.Range("AG40").Value = 8
$ws.Range("AI40").Value = 17

# Row 43
$ws.Range("G43").Value = 1.87
$ws.Range("J43").Value = 2.47
$ws.Range("K43").Value = 2
$ws.Range("L43").Value = 4.45
$ws.Range("Q43").Value = 1.85
$ws.Range("R43").Value = 1.75
$ws.Range("W43").Value = 1.4
$ws.Range("X43").Value = 2.4
$ws.Range("AA43").Value = 6.1
$ws.Range("AB43").Value = 7.7
$ws.Range("AC43").Value = 6.9
$ws.Range("AD43").Value = 13.5
$ws.Range("AF43").Value = 19.5
$ws.Range("AG43").Value = 9.25
$ws.Range("AH43").Value = 5.4
$ws.Range("AJ43").Value = 45
$ws.Range("AK43").Value = 300
$ws.Range("AL43").Value = 9.5
$ws.Range("AM43").Value = 18.5
$ws.Range("AN43").Value = 10.75
$ws.Range("AO43").Value = 50
$ws.Range("AP43").Value = 30
$ws.Range("AQ43").Value = 30

# Row 46
$ws.Range("G46").Value = 1.91
$ws.Range("I46").Value = 3.9
$ws.Range("J46").Value = 2.75
$ws.Range("L46").Value = 5
$ws.Range("W46").Value = 1.62
$ws.Range("X46").Value = 2.2
$ws.Range("Y46").Value = 2.25
$ws.Range("Z46").Value = 1.57
$ws.Range("AB46").Value = 8
$ws.Range("AC46").Value = 9.5
$ws.Range("AD46").Value = 17
$ws.Range("AL46").Value = 8.5
$ws.Range("AM46").Value = 19
$ws.Range("AN46").Value = 15
$ws.Range("AO46").Value = 51
$ws.Range("AR46").Value = 2
$ws.Range("AS46").Value = 1.8

# Row 48
$ws.Range("W48").Value = 1.29
$ws.Range("X48").Value = 3.5
$ws.Range("AH48").Value = 10

# Row 49
$ws.Range("L49").Value = 3.4
$ws.Range("AB49").Value = 12
$ws.Range("AL49").Value = 10
$ws.Range("AO49").Value = 29

# Row 51
$ws.Range("J51").Value = 2.75
$ws.Range("L51").Value = 3.25
$ws.Range("AE51").Value = 15
$ws.Range("AM51").Value = 19

# Row 52
$ws.Range("G52").Value = 2.25
$ws.Range("I52").Value = 2.8
$ws.Range("Y52").Value = 1.53
$ws.Range("Z52").Value = 2.38
$ws.Range("AA52").Value = 11
$ws.Range("AB52").Value = 13
$ws.Range("AC52").Value = 9.5
$ws.Range("AI52").Value = 12
$ws.Range("AQ52").Value = 23

# Row 53
$ws.Range("H53").Value = 3.2
$ws.Range("K53").Value = 2
$ws.Range("Q53").Value = 2.35
$ws.Range("R53").Value = 1.57
$ws.Range("Y53").Value = 2.1
$ws.Range("Z53").Value = 1.67
$ws.Range("AQ53").Value = 51

# Row 54
$ws.Range("Q54").Value = 2.4
$ws.Range("R54").Value = 1.53
$ws.Range("U54").Value = 4.5
$ws.Range("V54").Value = 1.18
$ws.Range("Y54").Value = 2
$ws.Range("Z54").Value = 1.73
$ws.Range("AF54").Value = 34
$ws.Range("AG54").Value = 7.5
$ws.Range("AH54").Value = 6
$ws.Range("AI54").Value = 17
$ws.Range("AL54").Value = 8
$ws.Range("AR54").Value = 1.8
$ws.Range("AS54").Value = 2

# Row 56
$ws.Range("G56").Value = 3.3
$ws.Range("H56").Value = 3.2
$ws.Range("I56").Value = 2.2
$ws.Range("J56").Value = 3.75
$ws.Range("K56").Value = 2.1
$ws.Range("L56").Value = 2.88
$ws.Range("M56").Value = 1.06
$ws.Range("N56").Value = 10
$ws.Range("AB56").Value = 17
$ws.Range("AC56").Value = 12
$ws.Range("AD56").Value = 34
$ws.Range("AE56").Value = 26
$ws.Range("AH56").Value = 6
$ws.Range("AM56").Value = 11
$ws.Range("AO56").Value = 21

# Row 58
$ws.Range("M58").Value = 1.07
$ws.Range("N58").Value = 9
$ws.Range("O58").Value = 1.36
$ws.Range("P58").Value = 3
$ws.Range("Q58").Value = 2.15
$ws.Range("R58").Value = 1.67

# Row 59
$ws.Range("J59").Value = 8
$ws.Range("M59").Value = 1.07
$ws.Range("N59").Value = 8.5
$ws.Range("Y59").Value = 2.38
$ws.Range("Z59").Value = 1.53
$ws.Range("AE59").Value = 67
$ws.Range("AG59").Value = 8.5
$ws.Range("AH59").Value = 8.5
$ws.Range("AI59").Value = 26
$ws.Range("AJ59").Value = 101
$ws.Range("AP59").Value = 15
$ws.Range("AQ59").Value = 41

# Row 60
$ws.Range("K60").Value = 2.4
$ws.Range("M60").Value = 1.04
$ws.Range("N60").Value = 13
$ws.Range("AD60").Value = 10
$ws.Range("AL60").Value = 17
$ws.Range("AM60").Value = 34

# Row 61
$ws.Range("Q61").Value = 2.35
$ws.Range("R61").Value = 1.57
$ws.Range("U61").Value = 4.33
$ws.Range("V61").Value = 1.2

# Row 63
$ws.Range("G63").Value = 1.85
$ws.Range("H63").Value = 3.7
$ws.Range("J63").Value = 2.4
$ws.Range("K63").Value = 2.38
$ws.Range("Q63").Value = 1.57
$ws.Range("R63").Value = 2.35
$ws.Range("AI63").Value = 13
$ws.Range("AJ63").Value = 41

# Row 64
$ws.Range("AK64").Value = 700

# Row 65
$ws.Range("H65").Value = 3.2
$ws.Range("I65").Value = 6
$ws.Range("K65").Value = 1.95
$ws.Range("L65").Value = 6.5
$ws.Range("M65").Value = 1.13
$ws.Range("N65").Value = 6
$ws.Range("Q65").Value = 2.6
$ws.Range("R65").Value = 1.48
$ws.Range("U65").Value = 5.5
$ws.Range("V65").Value = 1.14
$ws.Range("W65").Value = 1.57
$ws.Range("X65").Value = 2.25
$ws.Range("Y65").Value = 2.5
$ws.Range("Z65").Value = 1.5
$ws.Range("AE65").Value = 19
$ws.Range("AG65").Value = 6
$ws.Range("AI65").Value = 23
$ws.Range("AJ65").Value = 101
$ws.Range("AN65").Value = 21
$ws.Range("AQ65").Value = 67
$ws.Range("AR65").Value = 2
$ws.Range("AS65").Value = 1.85

# Row 68
$ws.Range("G68").Value = 2.37
$ws.Range("I68").Value = 2.8
$ws.Range("J68").Value = 2.9
$ws.Range("K68").Value = 2.15
$ws.Range("L68").Value = 3.35
$ws.Range("P68").Value = 3.9
$ws.Range("U68").Value = 2.52
$ws.Range("W68").Value = 1.35
$ws.Range("X68").Value = 2.95
$ws.Range("AA68").Value = 10.75
$ws.Range("AB68").Value = 14.5
$ws.Range("AC68").Value = 9
$ws.Range("AD68").Value = 27
$ws.Range("AE68").Value = 17
$ws.Range("AF68").Value = 21
$ws.Range("AJ68").Value = 40
$ws.Range("AL68").Value = 11
$ws.Range("AM68").Value = 16
$ws.Range("AN68").Value = 10
$ws.Range("AP68").Value = 21
$ws.Range("AQ68").Value = 25

# Row 69
$ws.Range("G69").Value = 4.75
$ws.Range("I69").Value = 1.7
$ws.Range("J69").Value = 4.75
$ws.Range("AB69").Value = 26
$ws.Range("AF69").Value = 41
$ws.Range("AO69").Value = 13
$ws.Range("AQ69").Value = 23

# Row 71
$ws.Range("G71").Value = 2.18
$ws.Range("H71").Value = 3.35
$ws.Range("I71").Value = 3.1
$ws.Range("J71").Value = 2.77
$ws.Range("K71").Value = 2.1
$ws.Range("L71").Value = 3.6
$ws.Range("M71").Value = 1.06
$ws.Range("N71").Value = 7.3
$ws.Range("O71").Value = 1.29
$ws.Range("P71").Value = 3.25
$ws.Range("Q71").Value = 1.87
$ws.Range("R71").Value = 1.83
$ws.Range("U71").Value = 3.05
$ws.Range("V71").Value = 1.33
$ws.Range("W71").Value = 1.4
$ws.Range("X71").Value = 2.72
$ws.Range("Y71").Value = 1.72
$ws.Range("Z71").Value = 2.02
$ws.Range("AA71").Value = 7.9
$ws.Range("AB71").Value = 10.75
$ws.Range("AC71").Value = 8.75
$ws.Range("AD71").Value = 21
$ws.Range("AE71").Value = 17.5
$ws.Range("AF71").Value = 27
$ws.Range("AG71").Value = 7.3
$ws.Range("AH71").Value = 6.4
$ws.Range("AI71").Value = 13.5
$ws.Range("AJ71").Value = 60
$ws.Range("AK71").Value = 450
$ws.Range("AL71").Value = 10
$ws.Range("AM71").Value = 16.5
$ws.Range("AN71").Value = 10.75
$ws.Range("AQ71").Value = 32

# Row 72
$ws.Range("G72").Value = 2.45
$ws.Range("H72").Value = 3.65
$ws.Range("I72").Value = 2.52
$ws.Range("J72").Value = 2.92
$ws.Range("K72").Value = 2.3
$ws.Range("L72").Value = 3
$ws.Range("M72").Value = 1.04
$ws.Range("N72").Value = 8.75
$ws.Range("P72").Value = 4.05
$ws.Range("R72").Value = 2.18
$ws.Range("X72").Value = 3.15
$ws.Range("Y72").Value = 1.53
$ws.Range("AA72").Value = 11
$ws.Range("AB72").Value = 14
$ws.Range("AC72").Value = 9.5
$ws.Range("AD72").Value = 26
$ws.Range("AE72").Value = 17.5
$ws.Range("AF72").Value = 23
$ws.Range("AG72").Value = 8.75
$ws.Range("AH72").Value = 7.3
$ws.Range("AI72").Value = 12
$ws.Range("AJ72").Value = 40
$ws.Range("AL72").Value = 11.25
$ws.Range("AM72").Value = 14.5
$ws.Range("AN72").Value = 9.75
$ws.Range("AO72").Value = 28
$ws.Range("AP72").Value = 18.5
$ws.Range("AQ72").Value = 23

# Row 73
$ws.Range("M73").Value = 1.13
$ws.Range("N73").Value = 6
$ws.Range("O73").Value = 1.67
$ws.Range("P73").Value = 2.1
$ws.Range("Q73").Value = 3.1
$ws.Range("R73").Value = 1.36

